# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2404   (the "before" format version)
#   *_new -> *_FV2410   (the "after" format version)
# Then (re)build the worksheet table over the data range so the table's
# column headers pick up the renamed values, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) cells ---------------------------------
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2404")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2410")
    }
}

# --- 2. Turn the data range into a real table (ListObject) ------------------
$lastRow = $usedRange.Rows.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
